$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "65.294.68"
$ws.Range("E2").Value2 = "  -0.95%  "
$ws.Range("D3").Value2 = "2.932.07"
$ws.Range("E3").Value2 = "  -2.79%  "
$ws.Range("E4").Value2 = "  -0.07%  "
$ws.Range("D5").Value2 = "'568.62"
$ws.Range("E5").Value2 = "  -2.95%  "
$ws.Range("D6").Value2 = "'158.38"
$ws.Range("E6").Value2 = "  +1.21%  "
$ws.Range("E7").Value2 = "  +0.11%  "
$ws.Range("D8").Value2 = "'0.516"
$ws.Range("E8").Value2 = "  -0.40%  "
$ws.Range("D9").Value2 = "2.924.93"
$ws.Range("E9").Value2 = "  -2.90%  "
$ws.Range("D10").Value2 = "'6.72"
$ws.Range("E10").Value2 = "  -4.95%  "
$ws.Range("E11").Value2 = "  -4.41%  "
$ws.Range("D12").Value2 = "'0.456"
$ws.Range("E12").Value2 = "  +0.49%  "
$ws.Range("E13").Value2 = "  -2.72%  "
$ws.Range("D14").Value2 = "'34.21"
$ws.Range("E14").Value2 = "  -0.92%  "
$ws.Range("E15").Value2 = "  -0.80%  "
$ws.Range("D16").Value2 = "65.331.86"
$ws.Range("E16").Value2 = "  -0.90%  "
$ws.Range("D17").Value2 = "3.420.61"
$ws.Range("E17").Value2 = "  -2.72%  "
$ws.Range("D18").Value2 = "'6.91"
$ws.Range("E18").Value2 = "  -0.84%  "
$ws.Range("D19").Value2 = "2.931.24"
$ws.Range("E19").Value2 = "  -2.78%  "
$ws.Range("D20").Value2 = "'15.38"
$ws.Range("E20").Value2 = "  +10.98%  "
$ws.Range("D21").Value2 = "'444.56"
$ws.Range("E21").Value2 = "  -4.22%  "
$ws.Range("D22").Value2 = "'0.686"
$ws.Range("E22").Value2 = "  +0.14%  "
$ws.Range("E23").Value2 = "  -2.14%  "
$ws.Range("D24").Value2 = "'82.10"
$ws.Range("E24").Value2 = "  -0.25%  "
$ws.Range("E25").Value2 = "  -1.86%  "
$ws.Range("D26").Value2 = "'12.08"
$ws.Range("E26").Value2 = "  -3.63%  "
$ws.Range("D27").Value2 = "'10.10"
$ws.Range("E27").Value2 = "  -6.22%  "
$ws.Range("E28").Value2 = "  +0.04%  "
$ws.Range("D29").Value2 = "'8.05"
$ws.Range("E29").Value2 = "  +1.26%  "
$ws.Range("D30").Value2 = "'2.39"
$ws.Range("E30").Value2 = "  -0.90%  "
$ws.Range("E31").Value2 = "  -1.89%  "
$ws.Range("E32").Value2 = "  -4.90%  "
$ws.Range("E33").Value2 = "  -0.08%  "
$ws.Range("E34").Value2 = "  -0.93%  "
$ws.Range("D35").Value2 = "'1.00"
$ws.Range("D36").Value2 = "'0.972"
$ws.Range("E36").Value2 = "  -2.52%  "
$ws.Range("D37").Value2 = "'5.72"
$ws.Range("E37").Value2 = "  -1.72%  "
$ws.Range("D38").Value2 = "'49.79"
$ws.Range("E38").Value2 = "  +0.99%  "
$ws.Range("D39").Value2 = "'44.85"
$ws.Range("E39").Value2 = "  -0.11%  "
$ws.Range("E40").Value2 = "  -9.72%  "
$ws.Range("B41").Value2 = "Kaspa"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value2 = "'0.120"
$ws.Range("E41").Value2 = "  -2.54%  "
$ws.Range("B42").Value2 = "TheGraph"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").Value2 = "'0.299"
$ws.Range("E42").Value2 = "  -0.34%  "
$ws.Range("D43").Value2 = "'2.83"
$ws.Range("E43").Value2 = "  -6.53%  "
$ws.Range("D44").Value2 = "'8.45"
$ws.Range("E44").Value2 = "  -0.60%  "
$ws.Range("D45").Value2 = "'384.27"
$ws.Range("E45").Value2 = "  -3.34%  "
$ws.Range("E46").Value2 = "  -1.16%  "
$ws.Range("D47").Value2 = "2.700.32"
$ws.Range("E47").Value2 = "  -3.76%  "
$ws.Range("D48").Value2 = "'133.25"
$ws.Range("E48").Value2 = "  -0.52%  "
$ws.Range("E50").Value2 = "  +4.01%  "
$ws.Range("D51").Value2 = "'23.34"
$ws.Range("E51").Value2 = "  -1.91%  "

# Reset style on cells that required a quote-prefix to stay text, so no stray style lingers
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").Style = "Normal"
